# "dua len bai tap" - add the two newly-published exercises (s11_vong lap 2)
# and their source links to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (link) for the two rows that already hold their exercise title
# in column A (A37 = "Sinh bang cuu chuong", A38 = "Hien thi cac so nguyen
# to dau tien") but were still missing their corresponding source link.
$ws.Range("B37").Value = 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%5D%20Sinh%20b%E1%BA%A3ng%20c%E1%BB%ADu%20ch%C6%B0%C6%A1ng.html'
$ws.Range("B38").Value = 'https://github.com/thanhviet05x1d/c0523g1_nguyen_thanh_viet_module1/blob/main/s11_vong%20lap%202/%5BB%C3%A0i%20t%E1%BA%ADp%5D%20Hi%E1%BB%83n%20th%E1%BB%8B%20c%C3%A1c%20s%E1%BB%91%20nguy%C3%AAn%20t%E1%BB%91%20%C4%91%E1%BA%A7u%20ti%C3%AAn.html'

# Reproduce the author's final on-screen state: scrolled up a bit with the
# two freshly-filled link cells selected.
$ws.Range("A23").Select()
$ws.Range("B37:B38").Select()
